$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# 1. Rename the original sheet and refresh its intro text (version bump).
$ws1.Name = "Self-assessment checklist"
$ws1.Range("A1").Value = "Onderstaande checklist kan gebruikt worden voor het uitvoeren van een assessment tegen de Kwaliteitsaanpak ICTU Software Realisatie versie 2.0-build.0, 14-08-2019."

# 2. Duplicate sheet1 (keeps identical column widths / row heights / styles)
#    right after it, then turn the duplicate into the new action-list sheet.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Self-assessment verbeteracties"

# Drop everything below the header rows - only the title + header row remain.
$ws2.Rows("3:76").Delete()

# Strip formatting/validation that only makes sense on the checklist sheet.
$ws2.Cells.FormatConditions.Delete()
$ws2.Cells.Validation.Delete()
$excel.ActiveWindow.FreezePanes = $false

for ($i = $ws2.Comments.Count; $i -ge 1; $i--) {
    $ws2.Comments.Item($i).Delete()
}

# New title + header row content for the action list sheet.
$ws2.Range("A1").Value = "Onderstaande actielijst kan gebruikt worden om acties n.a.v. de self-assessment bij te houden."
$ws2.Range("A2").Value = "Datum"
$ws2.Range("B2").Value = "Actie"
$ws2.Range("C2").Value = "Status"
$ws2.Range("D2").Value = "Toelichting"

$ws1.Select()
